# Updated cryptos list on Mon Jun 17 09:34:05 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a string value to a cell while preventing Excel's COM
# layer from auto-coercing numeric-looking text (e.g. "602.83") into a
# real number. We briefly force a text number-format, assign the value,
# then restore the cell style so no stray formatting is left behind.
function Set-TextValue($rangeAddr, $val) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "66.176.04"
$ws.Range("E2").Value = "  -0.33%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.545.25"
$ws.Range("E3").Value = "  -0.32%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
Set-TextValue "D5" "602.83"

# Row 6 - Solana
Set-TextValue "D6" "146.68"
$ws.Range("E6").Value = "  +1.48%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.547.57"
$ws.Range("E7").Value = "  -0.24%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.18%  "

# Row 9 - XRP
Set-TextValue "D9" "0.502"
$ws.Range("E9").Value = "  +2.21%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.134"
$ws.Range("E10").Value = "  -1.78%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  -0.70%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -1.06%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "4.143.06"
$ws.Range("E13").Value = "  -0.41%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  -2.40%  "

# Row 15 - Avalanche
Set-TextValue "D15" "29.06"
$ws.Range("E15").Value = "  -3.77%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.534.80"
$ws.Range("E16").Value = "  -0.50%  "

# Row 17 - TRON
$ws.Range("E17").Value = "  +1.50%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "66.129.27"
$ws.Range("E18").Value = "  -0.51%  "

# Row 19 - Uniswap
$ws.Range("E19").Value = "  -4.14%  "

# Row 20 - Polkadot
Set-TextValue "D20" "6.25"
$ws.Range("E20").Value = "  +1.00%  "

# Row 21 - Chainlink
Set-TextValue "D21" "14.66"
$ws.Range("E21").Value = "  -1.36%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "419.78"
$ws.Range("E22").Value = "  -2.54%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.603"
$ws.Range("E23").Value = "  -1.20%  "

# Row 24 - Litecoin
Set-TextValue "D24" "77.87"

# Row 25 - WrappedeETH
Set-TextValue "D25" "3.681.69"
$ws.Range("E25").Value = "  -0.42%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.06%  "

# Row 27 - PEPE
$ws.Range("E27").Value = "  -2.88%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue "D28" "9.14"
$ws.Range("E28").Value = "  -0.15%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  -1.58%  "

# Row 30 - RenderToken
Set-TextValue "D30" "7.83"
$ws.Range("E30").Value = "  -1.71%  "

# Row 31 - Binance-PegBSC-USD
$ws.Range("E31").Value = "  +0.13%  "

# Row 32 - RenzoRestakedETH
Set-TextValue "D32" "3.539.25"
$ws.Range("E32").Value = "  -0.33%  "

# Row 33 - Kaspa
Set-TextValue "D33" "0.155"
$ws.Range("E33").Value = "  +0.16%  "

# Row 34 - EthereumClassic
Set-TextValue "D34" "24.48"
$ws.Range("E34").Value = "  -3.88%  "

# Row 35 - USDe (unchanged)

# Row 36 - Aptos
Set-TextValue "D36" "7.60"
$ws.Range("E36").Value = "  -2.99%  "

# Row 37 - Fetch.AI
$ws.Range("E37").Value = "  -11.06%  "

# Row 38 - was NEARProtocol, now Monero
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D38" "174.18"
$ws.Range("E38").Value = "  -1.03%  "

# Row 39 - was Monero, now NEARProtocol
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D39" "5.28"
$ws.Range("E39").Value = "  -5.41%  "

# Row 40 - ImmutableX
$ws.Range("E40").Value = "  -7.39%  "

# Row 41 - Hedera
Set-TextValue "D41" "0.0827"
$ws.Range("E41").Value = "  -2.68%  "

# Row 42 - Filecoin
$ws.Range("E42").Value = "  -1.56%  "

# Row 43 - Mantle
Set-TextValue "D43" "0.862"
$ws.Range("E43").Value = "  -3.06%  "

# Row 44 - OKB
Set-TextValue "D44" "45.66"
$ws.Range("E44").Value = "  -0.76%  "

# Row 45 - Stacks
Set-TextValue "D45" "1.81"
$ws.Range("E45").Value = "  -6.05%  "

# Row 46 - FirstDigitalUSD
$ws.Range("E46").Value = "  +0.00%  "

# Row 47 - dogwifhat
Set-TextValue "D47" "2.42"
$ws.Range("E47").Value = "  -4.19%  "

# Row 48 - Cosmos
Set-TextValue "D48" "7.14"
$ws.Range("E48").Value = "  -0.10%  "

# Row 49 - ONDO
Set-TextValue "D49" "1.10"
$ws.Range("E49").Value = "  -7.31%  "

# Row 50 - EnergySwap
$ws.Range("E50").Value = "  -2.11%  "

# Row 51 - InjectiveProtocol
Set-TextValue "D51" "23.24"
$ws.Range("E51").Value = "  -7.51%  "
